$d = $word.ActiveDocument

# --- Body-text run-merge touch-ups ------------------------------------
# The source paragraph's text doesn't actually change, but three of the
# spell-checked words ("Morbi", "eros", "Sed") are no longer flagged, so
# their <w:proofErr> wrapping/run split disappears and they merge back
# into the neighbouring run. Doing a targeted Find/Replace over exactly
# the merged span reproduces that run layout without touching the other
# still-flagged words ("rhoncus", "facilisis", "odio").
$d.Content.Find.Execute(". Morbi id ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ". Morbi id ", 2) | Out-Null

$d.Content.Find.Execute(" eros, ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " eros, ", 2) | Out-Null

$d.Content.Find.Execute(". Sed ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ". Sed ", 2) | Out-Null

# --- Drop the stale _GoBack bookmark -----------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
